$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "Tue Sep 26 21:27:57 EDT 2023"
$ws.Range("B3").Value = "Tue Sep 26 21:28:12 EDT 2023"
$ws.Range("B4").Value = "Tue Sep 26 21:28:25 EDT 2023"
$ws.Range("B5").Value = "Tue Sep 26 21:28:39 EDT 2023"
